# Add a new "Tools to Install" slide at the end of the deck (slide 17),
# using the same "Title and Content" layout as the rest of the deck.

$p = $ppt.ActivePresentation

# Append a brand-new slide after the last one, using PpSlideLayout 2
# (ppLayoutText => Title + Content placeholder), matching slideLayout2.xml
# used by the other content slides in this deck.
$newIndex = $p.Slides.Count + 1
$s = $p.Slides.Add($newIndex, 2)

# --- Title placeholder: left blank (matches the source slide) ---
# (Do not touch it -- the freshly-added placeholder already starts empty.)

# --- Body / content placeholder ---
$body = $s.Shapes.Item(2)
$tf = $body.TextFrame
$tr = $tf.TextRange

$line1 = "Install Eclipse IDE Corresponding to Java 8 Version"
$line2 = "https://www.eclipse.org/downloads/packages/release/neon/3/eclipse-ide-java-ee-developers"
$line3 = "Git-https://git-scm.com/"
$line4 = "Tortoise Git -https://tortoisegit.org/download/"
$line5 = "Sighn up an account in git-> https://github.com/"

# Build the paragraphs up one at a time via InsertAfter (rather than a
# single multi-line Text= assignment) so every run keeps its lang="en-US".
$tr.Text = $line1
$tr.InsertAfter("`r" + $line2) | Out-Null
$tr.InsertAfter("`r" + $line3) | Out-Null
$tr.InsertAfter("`r" + $line4) | Out-Null
$tr.InsertAfter("`r" + $line5) | Out-Null
$tr.InsertAfter("`r") | Out-Null

# Paragraph 3: "Git" / "-https://git-scm.com/" -- re-stamp each sub-range
# with its own (unchanged) text so PowerPoint records them as distinct runs.
$para3 = $tr.Paragraphs(3, 1)
$p3r1 = $para3.Characters(1, 3)
$p3r1.Text = $p3r1.Text
$p3r2 = $para3.Characters(4, 21)
$p3r2.Text = $p3r2.Text

# Paragraph 4: "Tortoise " / "Git" / " -https://tortoisegit.org/download/"
$para4 = $tr.Paragraphs(4, 1)
$p4r1 = $para4.Characters(1, 9)
$p4r1.Text = $p4r1.Text
$p4r2 = $para4.Characters(10, 3)
$p4r2.Text = $p4r2.Text
$p4r3 = $para4.Characters(13, 35)
$p4r3.Text = $p4r3.Text

# Paragraph 5: "Sighn" / " up an account in " / "git" / "-> https://github.com/"
$para5 = $tr.Paragraphs(5, 1)
$p5r1 = $para5.Characters(1, 5)
$p5r1.Text = $p5r1.Text
$p5r2 = $para5.Characters(6, 18)
$p5r2.Text = $p5r2.Text
$p5r3 = $para5.Characters(24, 3)
$p5r3.Text = $p5r3.Text
$p5r4 = $para5.Characters(27, 22)
$p5r4.Text = $p5r4.Text

# Paragraph 6: blank trailing paragraph with no bullet.
$para6 = $tr.Paragraphs(6, 1)
$para6.ParagraphFormat.Bullet.Type = 0

# The amount of text overflows the placeholder, so PowerPoint shrinks it
# via "Shrink text on overflow" (normAutofit).
$tf.AutoSize = 2
